$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the old row 4: the header/week-grid block (old
# rows 4-11) and the footer notes (old rows 14-15) both shift down by 2,
# opening up new rows 4-5 and turning the old A1:E15 sheet into A1:E17.
$ws.Rows("4:5").Insert()

# The two custom column-width bands ("A"-style wide columns, "B"-style
# narrower day columns) each grow by two columns to keep covering the
# now-wider-looking grid: max 15 -> 17 and max 11 -> 13.
$ws.Range($ws.Cells.Item(1, 16), $ws.Cells.Item(1, 17)).ColumnWidth = 22.833333333333332
$ws.Range($ws.Cells.Item(1, 12), $ws.Cells.Item(1, 13)).ColumnWidth = 14.833333333333334

# The "Nombre d'équipes" line of the stats block (now at A17) goes from 4 to 0.
$null = $ws.Cells.Replace("Nombre d'équipes: 4", "Nombre d'équipes: 0", 2)

# Editing that multi-line cell makes the host auto-grow the row to fit the
# wrapped text; put row 17 back to the default unsized height.
$r = $ws.Rows(17)
$r.RowHeight = 15
$r.AutoFit()
